$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r = 2;  b = "01/01/2014"; c = 2343.71 },
    @{ r = 3;  b = "01/01/2015"; c = 2280.98 },
    @{ r = 4;  b = "01/01/2016"; c = 2243.6 },
    @{ r = 5;  b = "01/01/2017"; c = 2251.55 },
    @{ r = 6;  b = "01/01/2018"; c = 2273.95 },
    @{ r = 7;  b = "01/01/2019"; c = 2287.66 },
    @{ r = 8;  b = "01/01/2020"; c = 2242.61 },
    @{ r = 9;  b = "01/01/2021"; c = 2348 },
    @{ r = 10; b = "01/01/2022"; c = 2363.38 },
    @{ r = 11; b = "01/01/2023"; c = 2451.2 },
    @{ r = 12; b = "01/01/2014"; c = 1450.6 },
    @{ r = 13; b = "01/01/2015"; c = 1416.02 },
    @{ r = 14; b = "01/01/2016"; c = 1413.59 },
    @{ r = 15; b = "01/01/2017"; c = 1409.71 },
    @{ r = 16; b = "01/01/2018"; c = 1418.22 },
    @{ r = 17; b = "01/01/2019"; c = 1454.08 },
    @{ r = 18; b = "01/01/2020"; c = 1412.31 },
    @{ r = 19; b = "01/01/2021"; c = 1503.21 },
    @{ r = 20; b = "01/01/2022"; c = 1518.26 },
    @{ r = 21; b = "01/01/2023"; c = 1627.57 },
    @{ r = 22; b = "01/01/2014"; c = 1759.57 },
    @{ r = 23; b = "01/01/2015"; c = 1727.41 },
    @{ r = 24; b = "01/01/2016"; c = 1683.66 },
    @{ r = 25; b = "01/01/2017"; c = 1617.25 },
    @{ r = 26; b = "01/01/2018"; c = 1634.63 },
    @{ r = 27; b = "01/01/2019"; c = 1609.62 },
    @{ r = 28; b = "01/01/2020"; c = 1563.93 },
    @{ r = 29; b = "01/01/2021"; c = 1649.58 },
    @{ r = 30; b = "01/01/2022"; c = 1626.91 },
    @{ r = 31; b = "01/01/2023"; c = 1674.84 }
)

foreach ($row in $rows) {
    $bcell = $ws.Cells.Item($row.r, 2)
    $bcell.NumberFormat = "@"
    $bcell.Value = $row.b
    $bcell.ClearFormats()
    $ws.Cells.Item($row.r, 3).Value = $row.c
}
